# Generate Report for Handoff
# Updates status text, handoff/handback timestamps, and widens the
# status/date columns on the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Target rendered column width is ~17.216 "characters". The host only
# keeps whole-pixel column widths internally, so feed it the character
# value whose rounded pixel width lands closest to that target.
$newColWidth = 16.333333333333332

# --- Overview sheet ---------------------------------------------------
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-25 14:39:42"

$overview.Range("E:E").ColumnWidth = $newColWidth
$overview.Range("F:F").ColumnWidth = $newColWidth

# --- zh-cn sheet --------------------------------------------------------
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-25 14:39:37"

$zhcn.Range("C:C").ColumnWidth = $newColWidth

# --- de-de sheet --------------------------------------------------------
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-25 14:39:42"

$dede.Range("C:C").ColumnWidth = $newColWidth
